$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.806.59'
$ws.Range("E2").Value = '  +3.51%  '

$ws.Range("D3").Value = '3.700.60'
$ws.Range("E3").Value = '  +7.53%  '

$ws.Range("E4").Value = '  +0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '588.52'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.73%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '180.90'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.92%  '

$ws.Range("D7").Value = '3.695.03'
$ws.Range("E7").Value = '  +7.65%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.617'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.86%  '

$ws.Range("E9").Value = '  +0.09%  '

$ws.Range("E11").Value = '  +4.54%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '49.98'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.97%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000288'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.25%  '

$ws.Range("D14").Value = '4.310.54'
$ws.Range("E14").Value = '  +8.13%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '681.47'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.49%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '9.04'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.04%  '

$ws.Range("D17").Value = '3.703.32'
$ws.Range("E17").Value = '  +7.86%  '

$ws.Range("D18").Value = '71.915.90'

$ws.Range("E19").Value = '  +2.11%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.11'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.27%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.66'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.71%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.40'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +19.24%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.947'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.12%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '17.82'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.57%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '103.43'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.22%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.05'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.37%  '

$ws.Range("E27").Value = '  +5.22%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.34'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +6.27%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '35.61'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +5.73%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.31'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.82%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.36'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +6.69%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.20'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +11.34%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.30'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.15%  '

$ws.Range("B34").Value = 'Bittensor'
$ws.Range("C34").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '566.38'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.53%  '

$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.110'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.75%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '59.60'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.45%  '

$ws.Range("D37").Value = '3.766.47'

$ws.Range("E38").Value = '  -0.04%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.144'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.93%  '

$ws.Range("D40").Value = '0.0₃0778'
$ws.Range("E40").Value = '  +4.70%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '35.78'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.82%  '

$ws.Range("E42").Value = '  +5.33%  '

$ws.Range("E43").Value = '  +4.46%  '

$ws.Range("E44").Value = '  +9.28%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.353'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.95%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.92'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +8.93%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.37'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.55%  '

$ws.Range("E48").Value = '  +3.42%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.44'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.29%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.999'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.14%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '135.81'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.38%  '
